# Add a new "health_facility" choice option into the facility_types list
# (choices sheet), which alphabetizes/expands the admin region / facility
# type menu options, and switch the active tab over to the choices sheet.

$wb = $excel.ActiveWorkbook

$wsSurvey   = $wb.Worksheets.Item(1)   # survey
$wsChoices  = $wb.Worksheets.Item(2)   # choices
$wsSettings = $wb.Worksheets.Item(4)   # settings

# Insert a new row for the facility_types group (after "maternity", row 17)
# in the choices sheet. This pushes the existing blank separator row and
# every subsequent group down by one row.
$wsChoices.Rows("18:18").Insert()

$wsChoices.Range("A18").Value = "facility_types"
$wsChoices.Range("B18").Value = "health_facility"
$wsChoices.Range("C18").Value = "Health Facility"
$wsChoices.Range("D18").Value = "Facilidad de Salúd"

# Update the settings sheet's remembered selection without making it the
# active tab.
$wsSettings.Range("A48").Select()

# Make the choices sheet the active tab/selection.
$wsChoices.Activate()
$wsChoices.Range("D19").Select()
